$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the column headers in row 1: "_old" -> "_FV2404", "_new" -> "_FV2410"
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value()
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2404"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2410"
        }
    }
}

# Turn the data range into an Excel table ("Table1")
$rng = $ws.Range("A1:U69")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row
$ws.Activate()
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
